# Auto-generated Excel COM-interop script to apply the Cactuar_Profits value updates
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 292.7143
$ws.Range("I9").Value = 241.66667
$ws.Range("K9").Value = 241.66667
$ws.Range("M9").Value = -72.66667000000001
$ws.Range("H16").Value = 50000
$ws.Range("J16").Value = 50000
$ws.Range("L16").Value = 50000
$ws.Range("N16").Value = -50460
$ws.Range("H32").Value = 1773.0625
$ws.Range("I32").Value = 2024.5
$ws.Range("J32").Value = 1689.25
$ws.Range("K32").Value = 2024.5
$ws.Range("L32").Value = 1689.25
$ws.Range("M32").Value = -1698.5
$ws.Range("N32").Value = -2341.25
$ws.Range("H132").Value = 107527.92
$ws.Range("I132").Value = 250435
$ws.Range("K132").Value = 751305
$ws.Range("M132").Value = -748775
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 68055
$ws.Range("J140").Value = 67838.89
$ws.Range("L140").Value = 67838.89
$ws.Range("N140").Value = -78198.89

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3686.121
$ws.Range("I32").Value = 2008.9824
$ws.Range("J32").Value = 14308
$ws.Range("K32").Value = 2008.9824
$ws.Range("L32").Value = 14308
$ws.Range("M32").Value = -1721.9824
$ws.Range("N32").Value = -14882
$ws.Range("H61").Value = 13850.303
$ws.Range("I61").Value = 10231.857
$ws.Range("K61").Value = 10231.857
$ws.Range("M61").Value = -10019.857
$ws.Range("H74").Value = 12502835
$ws.Range("I74").Value = 19232130
$ws.Range("J74").Value = 5573.2856
$ws.Range("K74").Value = 19232130
$ws.Range("L74").Value = 5573.2856
$ws.Range("M74").Value = -19231256
$ws.Range("N74").Value = -7321.2856
$ws.Range("H77").Value = 12502835
$ws.Range("I77").Value = 19232130
$ws.Range("J77").Value = 5573.2856
$ws.Range("K77").Value = 96160650
$ws.Range("L77").Value = 27866.428
$ws.Range("M77").Value = -96156282
$ws.Range("N77").Value = -36602.428
$ws.Range("H132").Value = 19380.676
$ws.Range("I132").Value = 18283.229
$ws.Range("K132").Value = 54849.687
$ws.Range("M132").Value = -52319.687
$ws.Range("H136").Value = 13850.303
$ws.Range("I136").Value = 10231.857
$ws.Range("K136").Value = 30695.571
$ws.Range("M136").Value = -28145.571
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 100000
$ws.Range("J87").Value = 100000
$ws.Range("L87").Value = 100000
$ws.Range("N87").Value = -102496
$ws.Range("H90").Value = 100000
$ws.Range("J90").Value = 100000
$ws.Range("L90").Value = 300000
$ws.Range("N90").Value = -312480
$ws.Range("H107").Value = 111115864
$ws.Range("I107").Value = 111115864
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 111115864
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -111113944
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 88095.625
$ws.Range("J132").Value = 88095.625
$ws.Range("L132").Value = 88095.625
$ws.Range("N132").Value = -98215.625
$ws.Range("H134").Value = 1869.4546
$ws.Range("I134").Value = 1620.7
$ws.Range("K134").Value = 4862.1
$ws.Range("M134").Value = -2327.1
$ws.Range("H140").Value = 261387.5
$ws.Range("J140").Value = 261387.5
$ws.Range("L140").Value = 261387.5
$ws.Range("N140").Value = -271747.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13059.444
$ws.Range("I99").Value = 21859.5
$ws.Range("K99").Value = 21859.5
$ws.Range("M99").Value = -20361.5
$ws.Range("H126").Value = 13059.444
$ws.Range("I126").Value = 21859.5
$ws.Range("K126").Value = 65578.5
$ws.Range("M126").Value = -63108.5
$ws.Range("H132").Value = 37045090
$ws.Range("I132").Value = 40406676
$ws.Range("K132").Value = 121220028
$ws.Range("M132").Value = -121217498
$ws.Range("H134").Value = 1925.8182
$ws.Range("I134").Value = 1768.5
$ws.Range("K134").Value = 5305.5
$ws.Range("M134").Value = -2770.5
$ws.Range("H141").Value = 148831.78
$ws.Range("J141").Value = 148831.78
$ws.Range("L141").Value = 148831.78
$ws.Range("N141").Value = -159191.78

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 11500
$ws.Range("J54").Value = 11500
$ws.Range("L54").Value = 34500
$ws.Range("N54").Value = -35618
$ws.Range("H55").Value = 3477.6667
$ws.Range("I55").Value = 1500
$ws.Range("J55").Value = 3724.875
$ws.Range("K55").Value = 4500
$ws.Range("L55").Value = 11174.625
$ws.Range("M55").Value = -4323
$ws.Range("N55").Value = -11528.625
$ws.Range("H82").Value = 4071.4285
$ws.Range("J82").Value = 4666.6665
$ws.Range("L82").Value = 13999.9995
$ws.Range("N82").Value = -14811.9995
$ws.Range("H85").Value = 4071.4285
$ws.Range("J85").Value = 4666.6665
$ws.Range("L85").Value = 13999.9995
$ws.Range("N85").Value = -16807.9995
$ws.Range("H131").Value = 9093053
$ws.Range("J131").Value = 1543133.6
$ws.Range("L131").Value = 4629400.800000001
$ws.Range("N131").Value = -4639480.800000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H122").Value = 923450.3
$ws.Range("I122").Value = 1838234.4
$ws.Range("J122").Value = 8666.333000000001
$ws.Range("K122").Value = 5514703.199999999
$ws.Range("L122").Value = 25998.999
$ws.Range("M122").Value = -5512253.199999999
$ws.Range("N122").Value = -30898.999
$ws.Range("H126").Value = 3510.6428
$ws.Range("I126").Value = 2099.9565
$ws.Range("K126").Value = 6299.869499999999
$ws.Range("M126").Value = -3829.869499999999
$ws.Range("H127").Value = 38000
$ws.Range("J127").Value = 38000
$ws.Range("L127").Value = 38000
$ws.Range("N127").Value = -47920
$ws.Range("H132").Value = 3017.468
$ws.Range("I132").Value = 2958.8206
$ws.Range("J132").Value = 3303.375
$ws.Range("K132").Value = 8876.461800000001
$ws.Range("L132").Value = 9910.125
$ws.Range("M132").Value = -6346.461800000001
$ws.Range("N132").Value = -14970.125
$ws.Range("H140").Value = 38593.332
$ws.Range("J140").Value = 38593.332
$ws.Range("L140").Value = 38593.332
$ws.Range("N140").Value = -48953.332

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 55558556
$ws.Range("J40").Value = 83339336
$ws.Range("K40").Value = 55558556
$ws.Range("L40").Value = 83339336
$ws.Range("M40").Value = -55558420
$ws.Range("N40").Value = -83339608
$ws.Range("H54").Value = 59084
$ws.Range("J54").Value = 59084
$ws.Range("L54").Value = 59084
$ws.Range("N54").Value = -60372
$ws.Range("H55").Value = 282.8
$ws.Range("I55").Value = 343.33334
$ws.Range("K55").Value = 343.33334
$ws.Range("M55").Value = -170.33334
$ws.Range("H122").Value = 111125720
$ws.Range("I122").Value = 500003260
$ws.Range("J122").Value = 17850.715
$ws.Range("K122").Value = 1500009780
$ws.Range("L122").Value = 53552.145
$ws.Range("M122").Value = -1500007330
$ws.Range("N122").Value = -58452.145
$ws.Range("H132").Value = 4089.62
$ws.Range("I132").Value = 4090.6262
$ws.Range("K132").Value = 12271.8786
$ws.Range("M132").Value = -9741.8786
$ws.Range("H139").Value = 89500
$ws.Range("J139").Value = 89500
$ws.Range("L139").Value = 89500
$ws.Range("N139").Value = -99780

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 12511758
$ws.Range("I20").Value = 12511758
$ws.Range("K20").Value = 12511758
$ws.Range("M20").Value = -12511518
$ws.Range("H54").Value = 71842.71000000001
$ws.Range("J54").Value = 71842.71000000001
$ws.Range("L54").Value = 71842.71000000001
$ws.Range("N54").Value = -72882.71000000001
$ws.Range("H126").Value = 3465.7778
$ws.Range("I126").Value = 3400.25
$ws.Range("K126").Value = 10200.75
$ws.Range("M126").Value = -7730.75
$ws.Range("H136").Value = 5528.34
$ws.Range("I136").Value = 2533.843
$ws.Range("J136").Value = 8645.062
$ws.Range("K136").Value = 7601.529
$ws.Range("L136").Value = 25935.186
$ws.Range("M136").Value = -5051.529
$ws.Range("N136").Value = -31035.186
